$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.533.87'
$ws.Range("E2").Value = '  -0.42%  '

$ws.Range("D3").Value = '3.225.72'
$ws.Range("E3").Value = '  +0.45%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'602.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.13%  '

$ws.Range("D6").Value = "'156.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.10%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").Value = '3.226.96'
$ws.Range("E8").Value = '  +0.42%  '

$ws.Range("E9").Value = '  -1.88%  '

$ws.Range("E10").Value = '  +0.28%  '

$ws.Range("D11").Value = "'5.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.26%  '

$ws.Range("D12").Value = "'0.501"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.92%  '

$ws.Range("D13").Value = "'0.0000269"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.28%  '

$ws.Range("D14").Value = "'38.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.14%  '

$ws.Range("D15").Value = '3.760.12'

$ws.Range("D16").Value = '66.579.24'
$ws.Range("E16").Value = '  -0.39%  '

$ws.Range("D17").Value = '3.227.99'
$ws.Range("E17").Value = '  +0.51%  '

$ws.Range("D18").Value = "'7.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.79%  '

$ws.Range("E19").Value = '  +1.05%  '

$ws.Range("D20").Value = "'507.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.04%  '

$ws.Range("D21").Value = "'15.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.16%  '

$ws.Range("D22").Value = "'0.739"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.75%  '

$ws.Range("D23").Value = "'8.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.60%  '

$ws.Range("D24").Value = "'14.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.00%  '

$ws.Range("E25").Value = '  +0.56%  '

$ws.Range("D26").Value = "'0.165"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +82.08%  '

$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.07%  '

$ws.Range("D28").Value = "'3.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.23%  '

$ws.Range("D29").Value = "'9.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.48%  '

$ws.Range("D30").Value = "'2.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.81%  '

$ws.Range("D31").Value = "'2.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.33%  '

$ws.Range("D32").Value = "'6.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.39%  '

$ws.Range("D33").Value = "'28.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.48%  '

$ws.Range("E34").Value = '  +0.13%  '

$ws.Range("E35").Value = '  -6.38%  '

$ws.Range("D36").Value = "'6.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.51%  '

$ws.Range("D37").Value = "'55.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.66%  '

$ws.Range("E38").Value = '  +13.19%  '

$ws.Range("D39").Value = "'496.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.80%  '

$ws.Range("D40").Value = "'3.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.82%  '

$ws.Range("D41").Value = "'0.0420"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.48%  '

$ws.Range("D42").Value = "'0.128"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.49%  '

$ws.Range("D43").Value = "'8.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.68%  '

$ws.Range("D44").Value = "'0.294"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.77%  '

$ws.Range("D45").Value = '2.942.23'
$ws.Range("E45").Value = '  +1.58%  '

$ws.Range("E46").Value = '  -1.78%  '

$ws.Range("D47").Value = "'28.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.63%  '

$ws.Range("D48").Value = "'2.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.23%  '

$ws.Range("E49").Value = '  -0.06%  '

$ws.Range("E51").Value = '  -4.47%  '
